# "protecting when log out"
# Turns the old "ERP_User_Table" demo sheet into a hidden admin area
# (UserPermissionTable + UserPasswordTable), protects every visible
# worksheet + the workbook structure, and adds a new blank "Sheet4".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rebuild the old "ERP_User_Table" sheet into "UserPermissionTable"
# ---------------------------------------------------------------------
$permSheet = $wb.Worksheets.Item(1)
$permSheet.Cells.Clear()
$permSheet.Name = "UserPermissionTable"

$permHeaders = @("ID", "UserPasswordTable", "Sheet1", "Sheet2", "Sheet3", "UserPermissionTable", "Sheet4")
for ($c = 0; $c -lt $permHeaders.Count; $c++) {
    $permSheet.Cells.Item(1, $c + 1).Value = $permHeaders[$c]
}

$permSheet.Range("A2").Value = "admin"
$permSheet.Range("B2:G2").Value = "Writable"

$permSheet.Range("A3").Value = "guest"
$permSheet.Range("B3").Value = "Invisible"
$permSheet.Range("C3").Value = "ReadOnly"
$permSheet.Range("D3").Value = "ReadOnly"
$permSheet.Range("E3").Value = "ReadOnly"
$permSheet.Range("F3").Value = "Invisible"
$permSheet.Range("G3").Value = "ReadOnly"

$permSheet.Columns.Item(2).ColumnWidth = 29.256510416666668
$permSheet.Columns.Item(6).ColumnWidth = 17.799479166666668

$null = $permSheet.Range("B2:G1048576").Validation.Add(3, 1, 1, '"Writable,ReadOnly,Invisible"')

$null = $permSheet.Range("G7").Select()

# ---------------------------------------------------------------------
# 2) Insert a brand new "UserPasswordTable" sheet right after it
# ---------------------------------------------------------------------
$pwSheet = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$pwSheet.Name = "UserPasswordTable"

$pwSheet.Range("A1").Value = "ID"
$pwSheet.Range("B1").Value = "PASSWORD"
$pwSheet.Range("A2").Value = "admin"
$pwSheet.Range("B2").Value = "su2018"

$null = $pwSheet.Range("G8").Select()

# ---------------------------------------------------------------------
# 3) Protect the three pre-existing data sheets, and drop a value onto
#    "Sheet3"
# ---------------------------------------------------------------------
$sheetPassword = "su2018"

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Protect($sheetPassword)

$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet2.Protect($sheetPassword)

$sheet3 = $wb.Worksheets.Item("Sheet3")
$sheet3.Range("A1").Value = 123
$sheet3.Protect($sheetPassword)

# ---------------------------------------------------------------------
# 4) Append a brand new, empty "Sheet4" at the end and make it active
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet4 = $wb.Worksheets.Add($null, $lastSheet)
$sheet4.Name = "Sheet4"
$sheet4.PageSetup.PaperSize = 9
$sheet4.PageSetup.Orientation = 1
$sheet4.Protect($sheetPassword)
$sheet4.Activate()

# ---------------------------------------------------------------------
# 5) Hide the two admin sheets very deeply (not reachable from the UI)
# ---------------------------------------------------------------------
$permSheet.Visible = 2
$pwSheet.Visible = 2

# ---------------------------------------------------------------------
# 6) Lock the workbook structure so sheets can't be un-hidden / moved /
#    deleted / renamed without the password, and scroll the tab strip
#    so the first visible tab is "Sheet1"
# ---------------------------------------------------------------------
$excel.ActiveWindow.DisplayedFirstSheet = 3
$wb.Protect($sheetPassword)
